$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14912
$ws1.Range("F3").Value = 18631
$ws1.Range("F5").Value = 118
$ws1.Range("F7").Value = 221
$ws1.Range("F13").Value = 52
$ws1.Range("F14").Value = 116
$ws1.Range("F15").Value = 202
$ws1.Range("F16").Value = 55
$ws1.Range("F17").Value = 1426
$ws1.Range("F20").Value = 87
$ws1.Range("F22").Value = 7734
$ws1.Range("F24").Value = 22
$ws1.Range("F26").Value = 1225
$ws1.Range("F29").Value = 107
$ws1.Range("F30").Value = 66
$ws1.Range("F31").Value = 157
$ws1.Range("F33").Value = 263
$ws1.Range("F34").Value = 5338
$ws1.Range("F35").Value = 26

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14912
$ws4.Range("F3").Value = 18631
$ws4.Range("F5").Value = 118
$ws4.Range("F7").Value = 221
$ws4.Range("F13").Value = 52
$ws4.Range("F14").Value = 116
$ws4.Range("F15").Value = 202
$ws4.Range("F16").Value = 55
$ws4.Range("F17").Value = 1426
$ws4.Range("F21").Value = 87
$ws4.Range("F23").Value = 7734
$ws4.Range("F25").Value = 22
$ws4.Range("F27").Value = 1225
$ws4.Range("F32").Value = 107
$ws4.Range("F33").Value = 66
$ws4.Range("F34").Value = 157
$ws4.Range("F35").Value = 0
$ws4.Range("F36").Value = 263
$ws4.Range("F37").Value = 5338
$ws4.Range("F38").Value = 26

